# Generate Report for Handoff
#
# Inserts a new row for the source file
# "d63a0061-871f-4bb6-aff5-9a1504adccbc.md" into all three worksheets
# (Overview, zh-cn, de-de) immediately before the existing
# "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md" row, pushing that row (and its
# hyperlinks) down by one. The new row carries "Ready for handoff" status
# with handoff timestamp 2016-03-24 02:19:03.

$wb = $excel.ActiveWorkbook

function Replace-RowHyperlink {
    param(
        $ws,
        [string]$oldAddr,
        [string]$newAddr,
        [string]$target,
        [string]$display
    )
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $oldAddr) {
            $hl.Delete()
        }
    }
    $null = $ws.Hyperlinks.Add($ws.Range($newAddr), $target, "", "", $display)
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn), C (de-de),
# D (Latest Handoff Date). Old row 7 (d68ca6e4) moves to row 8; new row 7
# holds the d63a0061 entry.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows(7).Insert()

$wsOverview.Range("A7").Value2 = "d63a0061-871f-4bb6-aff5-9a1504adccbc.md"
$wsOverview.Range("B7").Value2 = "Ready for handoff"
$wsOverview.Range("C7").Value2 = "Ready for handoff"
$wsOverview.Range("D7").Value2 = "2016-03-24 02:19:03"

Replace-RowHyperlink -ws $wsOverview -oldAddr '$A$7' -newAddr "A8" `
    -target "https://github.com/OpenLocalizationTest/oltest/blob/5ba180692774f67357fa7c1130a929382a060d7a/e2e/d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md" `
    -display "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md"

$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/d63a0061871f4bb6aff59a1504adccbc000001/e2e/d63a0061-871f-4bb6-aff5-9a1504adccbc.md", `
    "", "", "d63a0061-871f-4bb6-aff5-9a1504adccbc.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime),
# H (Latest Handback DateTime), J (Handoff Reason).
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows(7).Insert()

$wsZh.Range("A7").Value2 = "d63a0061-871f-4bb6-aff5-9a1504adccbc.md"
$wsZh.Range("B7").Value2 = ".md"
$wsZh.Range("C7").Value2 = "Ready for handoff"
$wsZh.Range("D7").Value2 = "d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.zh-cn.xlf"
$wsZh.Range("E7").Value2 = "2016-03-24 02:18:53"
$wsZh.Range("H7").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("J7").Value2 = "Include"

Replace-RowHyperlink -ws $wsZh -oldAddr '$A$7' -newAddr "A8" `
    -target "https://github.com/OpenLocalizationTest/oltest/blob/5ba180692774f67357fa7c1130a929382a060d7a/e2e/d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md" `
    -display "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md"
Replace-RowHyperlink -ws $wsZh -oldAddr '$D$7' -newAddr "D8" `
    -target "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9d4a12baaa1edc6f8527419a1ff1f703a865d42/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d68ca6e4-a0cc-4175-8309-ddb4a376cf28.991d147a5b4d6a766dc6d5a14d8d9353aa9e654d.zh-cn.xlf" `
    -display "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.991d147a5b4d6a766dc6d5a14d8d9353aa9e654d.zh-cn.xlf"

$null = $wsZh.Hyperlinks.Add($wsZh.Range("A7"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/d63a0061871f4bb6aff59a1504adccbc000001/e2e/d63a0061-871f-4bb6-aff5-9a1504adccbc.md", `
    "", "", "d63a0061-871f-4bb6-aff5-9a1504adccbc.md")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("D7"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d63a0061871f4bb6aff59a1504adccbc000002/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.zh-cn.xlf", `
    "", "", "d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout, de-de XLIFF filenames/timestamps.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows(7).Insert()

$wsDe.Range("A7").Value2 = "d63a0061-871f-4bb6-aff5-9a1504adccbc.md"
$wsDe.Range("B7").Value2 = ".md"
$wsDe.Range("C7").Value2 = "Ready for handoff"
$wsDe.Range("D7").Value2 = "d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.de-de.xlf"
$wsDe.Range("E7").Value2 = "2016-03-24 02:19:03"
$wsDe.Range("H7").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("J7").Value2 = "Include"

Replace-RowHyperlink -ws $wsDe -oldAddr '$A$7' -newAddr "A8" `
    -target "https://github.com/OpenLocalizationTest/oltest/blob/5ba180692774f67357fa7c1130a929382a060d7a/e2e/d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md" `
    -display "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md"
Replace-RowHyperlink -ws $wsDe -oldAddr '$D$7' -newAddr "D8" `
    -target "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f811ade9f1cb03314d0fd4962d19fd1dd1a4b88e/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d68ca6e4-a0cc-4175-8309-ddb4a376cf28.991d147a5b4d6a766dc6d5a14d8d9353aa9e654d.de-de.xlf" `
    -display "d68ca6e4-a0cc-4175-8309-ddb4a376cf28.991d147a5b4d6a766dc6d5a14d8d9353aa9e654d.de-de.xlf"

$null = $wsDe.Hyperlinks.Add($wsDe.Range("A7"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/d63a0061871f4bb6aff59a1504adccbc000001/e2e/d63a0061-871f-4bb6-aff5-9a1504adccbc.md", `
    "", "", "d63a0061-871f-4bb6-aff5-9a1504adccbc.md")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("D7"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d63a0061871f4bb6aff59a1504adccbc000003/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.de-de.xlf", `
    "", "", "d63a0061-871f-4bb6-aff5-9a1504adccbc.7c04be0188bb06e3689654262266c326915c6cdf.de-de.xlf")
